# Updated cryptos list on Tue Apr  4 19:46:22 UTC 2023 with GitHub Actions
#
# Refreshes the Price (D) and Volume(1h) (E) columns with the latest scrape,
# and (for this run) corrects the EthereumClassic / Wrapped liquid staked
# Ether 2.0 rows, which had been swapped in the previous commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Helper: write a value as literal text, the same way Excel stores it when a
# user types a leading apostrophe in front of a numeric-looking entry -- this
# keeps price strings like "311.76" or "1.000" from being reinterpreted as
# numbers (and losing trailing zeros) while leaving already-non-numeric text
# (e.g. "28.188.19") untouched.
function Set-TextValue($range, [string]$text) {
    if ($text -match '^[+-]?\d+(\.\d+)?$') {
        $range.Value = "'" + $text
    } else {
        $range.Value = $text
    }
}

Set-TextValue $ws.Range("D2") "28.188.19"
$ws.Range("E2").Value = "  +0.13%  "
Set-TextValue $ws.Range("D3") "1.871.23"
$ws.Range("E3").Value = "  +2.36%  "
$ws.Range("E4").Value = "  +0.32%  "
Set-TextValue $ws.Range("D5") "311.76"
$ws.Range("E5").Value = "  -0.13%  "
$ws.Range("E6").Value = "  +0.09%  "
Set-TextValue $ws.Range("D7") "0.5057"
$ws.Range("E7").Value = "  -1.42%  "
Set-TextValue $ws.Range("D8") "0.3927"
$ws.Range("E8").Value = "  -0.08%  "
Set-TextValue $ws.Range("D9") "0.09688"
$ws.Range("E9").Value = "  -5.37%  "
Set-TextValue $ws.Range("D11") "40.89"
$ws.Range("E11").Value = "  -0.32%  "
Set-TextValue $ws.Range("D12") "6.514"
$ws.Range("E12").Value = "  +0.24%  "
Set-TextValue $ws.Range("D13") "20.93"
$ws.Range("E13").Value = "  +0.17%  "
Set-TextValue $ws.Range("D14") "1.879.84"
$ws.Range("E14").Value = "  +3.00%  "
Set-TextValue $ws.Range("D15") "1.002"
$ws.Range("E15").Value = "  +0.38%  "
Set-TextValue $ws.Range("D16") "7.407"
$ws.Range("E16").Value = "  -0.03%  "
$ws.Range("E17").Value = "  -1.39%  "
Set-TextValue $ws.Range("D18") "92.85"
$ws.Range("E18").Value = "  -1.75%  "
Set-TextValue $ws.Range("D19") "0.06591"
$ws.Range("E19").Value = "  -0.32%  "
Set-TextValue $ws.Range("D20") "17.56"
$ws.Range("E20").Value = "  +1.02%  "
$ws.Range("E21").Value = "  +0.05%  "
Set-TextValue $ws.Range("D22") "6.166"
$ws.Range("E22").Value = "  +1.83%  "
Set-TextValue $ws.Range("D23") "28.231.60"
$ws.Range("E23").Value = "  +0.00%  "
$ws.Range("E24").Value = "  +1.26%  "
Set-TextValue $ws.Range("D25") "2.276"
$ws.Range("E25").Value = "  +1.40%  "
Set-TextValue $ws.Range("D26") "2.535"
$ws.Range("E26").Value = "  +2.88%  "
$ws.Range("B27").Value = "EthereumClassic"
$ws.Range("C27").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
Set-TextValue $ws.Range("D27") "21.25"
$ws.Range("E27").Value = "  +1.73%  "
$ws.Range("B28").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C28").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
Set-TextValue $ws.Range("D28") "2.074.59"
$ws.Range("E28").Value = "  +1.82%  "
Set-TextValue $ws.Range("D29") "158.25"
$ws.Range("E29").Value = "  -0.31%  "
Set-TextValue $ws.Range("D30") "127.50"
$ws.Range("E30").Value = "  -1.20%  "
Set-TextValue $ws.Range("D31") "0.1062"
$ws.Range("E31").Value = "  -2.83%  "
Set-TextValue $ws.Range("D32") "1.068"
$ws.Range("E32").Value = "  -0.56%  "
Set-TextValue $ws.Range("D33") "5.631"
$ws.Range("E33").Value = "  -0.45%  "
Set-TextValue $ws.Range("D34") "3.627"
$ws.Range("E34").Value = "  -0.12%  "
$ws.Range("E35").Value = "  +4.40%  "
Set-TextValue $ws.Range("D36") "0.06712"
$ws.Range("E36").Value = "  -3.69%  "
Set-TextValue $ws.Range("D37") "0.02393"
$ws.Range("E37").Value = "  +1.83%  "
Set-TextValue $ws.Range("D38") "0.2181"
$ws.Range("E38").Value = "  -0.04%  "
$ws.Range("E39").Value = "  -1.64%  "
Set-TextValue $ws.Range("D40") "0.6358"
$ws.Range("E40").Value = "  +1.13%  "
Set-TextValue $ws.Range("D41") "4.967"
$ws.Range("E41").Value = "  -1.27%  "
Set-TextValue $ws.Range("D42") "1.176"
$ws.Range("E42").Value = "  +1.47%  "
Set-TextValue $ws.Range("D43") "1.000"
$ws.Range("E43").Value = "  +0.07%  "
$ws.Range("E44").Value = "  +1.38%  "
$ws.Range("E45").Value = "  -0.13%  "
Set-TextValue $ws.Range("D46") "3.658"
Set-TextValue $ws.Range("D47") "1.258"
$ws.Range("E47").Value = "  -2.25%  "
Set-TextValue $ws.Range("D48") "124.06"
$ws.Range("E48").Value = "  -1.51%  "
Set-TextValue $ws.Range("D49") "1.988"
$ws.Range("E49").Value = "  -0.93%  "
Set-TextValue $ws.Range("D50") "1.196"
$ws.Range("E50").Value = "  +0.28%  "
Set-TextValue $ws.Range("D51") "0.06837"
$ws.Range("E51").Value = "  +0.72%  "
